# Workbook / worksheet handles.
# The edited sheet is "strategy_id-0", which is the active sheet (dimension
# A1:AS11 -> A1:AS12 in the diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row above the current row 4
# ("elasticity_gnrl_rate_occupancy_to_gdppc"), pushing it and every row
# below it (frac_gnrl_eating_red_meat, limit_gnrl_annual_emissions_mt_ch4,
# limit_gnrl_annual_emissions_mt_co2, limit_gnrl_annual_emissions_mt_n2o,
# occrateinit_gnrl_occupancy, population_gnrl_rural, population_gnrl_urban)
# down by one row.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the
# "climate_change_factor_gnrl_hydropower_availability" variable.
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5
$ws.Range("J4:AS4").Value = 1
